$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # quality_comparison
$ws2 = $wb.Worksheets.Item(2)   # computational_comparison

# --- quality_comparison sheet ---

# Build the "top+bottom" border style on C1 (matches target borderId=4)
$c1 = $ws1.Range("C1")
$c1.ClearFormats()
$c1.Borders.Item(8).LineStyle = 1   # xlEdgeTop
$c1.Borders.Item(9).LineStyle = 1   # xlEdgeBottom

# Build the "top+right+bottom" border style on D1 (matches target borderId=5)
$d1 = $ws1.Range("D1")
$d1.ClearFormats()
$d1.Borders.Item(9).LineStyle = 1    # xlEdgeBottom
$d1.Borders.Item(10).LineStyle = 1   # xlEdgeRight
$d1.Borders.Item(8).LineStyle = 1    # xlEdgeTop

# anonymize column header
$ws1.Range("C2").Value = "approach"

# --- computational_comparison sheet ---

# Re-use the already-resolved styles from sheet1 via copy/paste-format so no
# extra intermediate style entries get created in the shared style table.
$c1.Copy()
$ws2.Range("C1").PasteSpecial(-4122)   # xlPasteFormats
$ws2.Range("F1").PasteSpecial(-4122)

$d1.Copy()
$ws2.Range("D1").PasteSpecial(-4122)
$ws2.Range("G1").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# anonymize column headers
$ws2.Range("C2").Value = "approach"
$ws2.Range("F2").Value = "approach"

# remove stray empty inline-string cell
$ws2.Range("G5").ClearContents()

Write-Host "edit.ps1 applied"
